# Update "想去人数" (number of people interested) counts for a couple of
# events that appear on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 648
$ws1.Range("F4").Value = 1483
$ws1.Range("F5").Value = 690

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 648
$ws4.Range("F4").Value = 1483
$ws4.Range("F6").Value = 690
